$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the time-range labels in column C (rows 7-12): shift from 19:xx to 21:xx
$ws.Range("C7").Value = "21:15-21:20"
$ws.Range("C8").Value = "21:20-21:25"
$ws.Range("C9").Value = "21:25-21:30"
$ws.Range("C10").Value = "21:30-21:35"
$ws.Range("C11").Value = "21:35-21:40"
$ws.Range("C12").Value = "21:40-21:45"

# Update the active selection to match the recorded cursor position
$ws.Range("C16").Select()
